$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.25226
$ws.Range("H2").Value = 3.75678
$ws.Range("I2").Value = 0.05753803679167191
$ws.Range("J2").Value = 0.06158044274193954
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7990306666666668
$ws.Range("N2").Value = 2.397092
$ws.Range("O2").Value = 0.3065513402634041
$ws.Range("P2").Value = 0.3495008831958056
$ws.Range("Q2").Value = 1.00059414264
$ws.Range("R2").Value = 9.005347283760001
$ws.Range("S2").Value = 0.01763836229461208
$ws.Range("T2").Value = 0.02152241912589661
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.25226
$ws.Range("H3").Value = 3.75678
$ws.Range("I3").Value = 0.05753803679167191
$ws.Range("J3").Value = 0.06158044274193954
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8465543333333333
$ws.Range("N3").Value = 2.539663
$ws.Range("O3").Value = 0.3247839867920704
$ws.Range("P3").Value = 0.3702880246230471
$ws.Range("Q3").Value = 1.06010612946
$ws.Range("R3").Value = 9.54095516514
$ws.Range("S3").Value = 0.01868743298138803
$ws.Range("T3").Value = 0.02280250049832545
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.25226
$ws.Range("H4").Value = 3.75678
$ws.Range("I4").Value = 0.05753803679167191
$ws.Range("J4").Value = 0.06158044274193954
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 0.9609300000000001
$ws.Range("N4").Value = 1.92186
$ws.Range("O4").Value = 0.3686646729445256
$ws.Range("P4").Value = 0.2802110921811474
$ws.Range("Q4").Value = 1.2033342018
$ws.Range("R4").Value = 7.2200052108
$ws.Range("S4").Value = 0.0212122415156718
$ws.Range("T4").Value = 0.01725552311771749
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.48419333333333
$ws.Range("H5").Value = 37.45258
$ws.Range("I5").Value = 0.5736156831070852
$ws.Range("J5").Value = 0.6139157624955174
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7990306666666668
$ws.Range("N5").Value = 2.397092
$ws.Range("O5").Value = 0.3065513402634041
$ws.Range("P5").Value = 0.3495008831958056
$ws.Range("Q5").Value = 9.975253321928889
$ws.Range("R5").Value = 89.77727989736
$ws.Range("S5").Value = 0.1758426564525851
$ws.Range("T5").Value = 0.2145641012000097
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.48419333333333
$ws.Range("H6").Value = 37.45258
$ws.Range("I6").Value = 0.5736156831070852
$ws.Range("J6").Value = 0.6139157624955174
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8465543333333333
$ws.Range("N6").Value = 2.539663
$ws.Range("O6").Value = 0.3247839867920704
$ws.Range("P6").Value = 0.3702880246230471
$ws.Range("Q6").Value = 10.56854796450444
$ws.Range("R6").Value = 95.11693168053999
$ws.Range("S6").Value = 0.186301188445976
$ws.Range("T6").Value = 0.2273256549794168
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.48419333333333
$ws.Range("H7").Value = 37.45258
$ws.Range("I7").Value = 0.5736156831070852
$ws.Range("J7").Value = 0.6139157624955174
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.9609300000000001
$ws.Range("N7").Value = 1.92186
$ws.Range("O7").Value = 0.3686646729445256
$ws.Range("P7").Value = 0.2802110921811474
$ws.Range("Q7").Value = 11.9964358998
$ws.Range("R7").Value = 71.9786153988
$ws.Range("S7").Value = 0.2114718382085242
$ws.Range("T7").Value = 0.1720260063160908
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.405587666666666
$ws.Range("H8").Value = 4.216762999999999
$ws.Range("I8").Value = 0.06458303777058033
$ws.Range("J8").Value = 0.06912039897939969
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7990306666666668
$ws.Range("N8").Value = 2.397092
$ws.Range("O8").Value = 0.3065513402634041
$ws.Range("P8").Value = 0.3495008831958056
$ws.Range("Q8").Value = 1.123107650355111
$ws.Range("R8").Value = 10.107968853196
$ws.Range("S8").Value = 0.01979801678685345
$ws.Range("T8").Value = 0.02415764049014665
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.405587666666666
$ws.Range("H9").Value = 4.216762999999999
$ws.Range("I9").Value = 0.06458303777058033
$ws.Range("J9").Value = 0.06912039897939969
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8465543333333333
$ws.Range("N9").Value = 2.539663
$ws.Range("O9").Value = 0.3247839867920704
$ws.Range("P9").Value = 0.3702880246230471
$ws.Range("Q9").Value = 1.189906330096555
$ws.Range("R9").Value = 10.709156970869
$ws.Range("S9").Value = 0.02097553648627195
$ws.Range("T9").Value = 0.02559445599923879
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.405587666666666
$ws.Range("H10").Value = 4.216762999999999
$ws.Range("I10").Value = 0.06458303777058033
$ws.Range("J10").Value = 0.06912039897939969
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 0.9609300000000001
$ws.Range("N10").Value = 1.92186
$ws.Range("O10").Value = 0.3686646729445256
$ws.Range("P10").Value = 0.2802110921811474
$ws.Range("Q10").Value = 1.35067135653
$ws.Range("R10").Value = 8.104028139179999
$ws.Range("S10").Value = 0.02380948449745494
$ws.Range("T10").Value = 0.01936830249001425
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.335940666666667
$ws.Range("H11").Value = 7.007822
$ws.Range("I11").Value = 0.1073302988371658
$ws.Range("J11").Value = 0.1148709217512615
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.7990306666666668
$ws.Range("N11").Value = 2.397092
$ws.Range("O11").Value = 0.3065513402634041
$ws.Range("P11").Value = 0.3495008831958056
$ws.Range("Q11").Value = 1.866488228180445
$ws.Range("R11").Value = 16.798394053624
$ws.Range("S11").Value = 0.03290224695940487
$ws.Range("T11").Value = 0.04014748860558218
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.335940666666667
$ws.Range("H12").Value = 7.007822
$ws.Range("I12").Value = 0.1073302988371658
$ws.Range("J12").Value = 0.1148709217512615
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.8465543333333333
$ws.Range("N12").Value = 2.539663
$ws.Range("O12").Value = 0.3247839867920704
$ws.Range("P12").Value = 0.3702880246230471
$ws.Range("Q12").Value = 1.977500693776222
$ws.Range("R12").Value = 17.797506243986
$ws.Range("S12").Value = 0.03485916235991903
$ws.Range("T12").Value = 0.04253532670190324
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.335940666666667
$ws.Range("H13").Value = 7.007822
$ws.Range("I13").Value = 0.1073302988371658
$ws.Range("J13").Value = 0.1148709217512615
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.9609300000000001
$ws.Range("N13").Value = 1.92186
$ws.Range("O13").Value = 0.3686646729445256
$ws.Range("P13").Value = 0.2802110921811474
$ws.Range("Q13").Value = 2.24467546482
$ws.Range("R13").Value = 13.46805278892
$ws.Range("S13").Value = 0.03956888951784194
$ws.Range("T13").Value = 0.03218810644377611
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 4.286056
$ws.Range("H14").Value = 8.572112000000001
$ws.Range("I14").Value = 0.1969329434934967
$ws.Range("J14").Value = 0.1405124740318818
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.7990306666666668
$ws.Range("N14").Value = 2.397092
$ws.Range("O14").Value = 0.3065513402634041
$ws.Range("P14").Value = 0.3495008831958056
$ws.Range("Q14").Value = 3.424690183050668
$ws.Range("R14").Value = 20.548141098304
$ws.Range("S14").Value = 0.06037005776994864
$ws.Range("T14").Value = 0.04910923377417038
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 4.286056
$ws.Range("H15").Value = 8.572112000000001
$ws.Range("I15").Value = 0.1969329434934967
$ws.Range("J15").Value = 0.1405124740318818
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.8465543333333333
$ws.Range("N15").Value = 2.539663
$ws.Range("O15").Value = 0.3247839867920704
$ws.Range("P15").Value = 0.3702880246230471
$ws.Range("Q15").Value = 3.628379279709333
$ws.Range("R15").Value = 21.770275678256
$ws.Range("S15").Value = 0.06396066651851537
$ws.Range("T15").Value = 0.05203008644416271
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 4.286056
$ws.Range("H16").Value = 8.572112000000001
$ws.Range("I16").Value = 0.1969329434934967
$ws.Range("J16").Value = 0.1405124740318818
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 0.9609300000000001
$ws.Range("N16").Value = 1.92186
$ws.Range("O16").Value = 0.3686646729445256
$ws.Range("P16").Value = 0.2802110921811474
$ws.Range("Q16").Value = 4.11859979208
$ws.Range("R16").Value = 16.47439916832
$ws.Range("S16").Value = 0.07260221920503269
$ws.Range("T16").Value = 0.03937315381354871